$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.034024354967587
$ws.Range("D2").Value = 1.036011125204094
$ws.Range("E2").Value = 1.04340624557467
$ws.Range("F2").Value = 1.054753321932304
$ws.Range("I2").Value = 1.035491254270328
$ws.Range("J2").Value = 1.039145866203091
$ws.Range("K2").Value = 1.038806088669119
$ws.Range("L2").Value = 1.046180201262493
$ws.Range("M2").Value = 1.057495664146019
$ws.Range("N2").Value = 1.016961110875613

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.035180448194092
$ws.Range("D3").Value = 1.036851503956108
$ws.Range("E3").Value = 1.044455552952765
$ws.Range("F3").Value = 1.055942122975875
$ws.Range("I3").Value = 1.03574971395546
$ws.Range("J3").Value = 1.039943872388502
$ws.Range("K3").Value = 1.039455898773381
$ws.Range("L3").Value = 1.047039914152897
$ws.Range("M3").Value = 1.058496811906301
$ws.Range("N3").Value = 1.017229738305801

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.035928238354883
$ws.Range("D4").Value = 1.037394718098104
$ws.Range("E4").Value = 1.045134629326813
$ws.Range("F4").Value = 1.056711536758896
$ws.Range("I4").Value = 1.035915057723964
$ws.Range("J4").Value = 1.040459439401277
$ws.Range("K4").Value = 1.039875160806197
$ws.Range("L4").Value = 1.047595710871174
$ws.Range("M4").Value = 1.059144228296668
$ws.Range("N4").Value = 1.017403177081639

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.03624254380014
$ws.Range("D5").Value = 1.037622949675225
$ws.Range("E5").Value = 1.045420138629737
$ws.Range("F5").Value = 1.057035042302833
$ws.Range("I5").Value = 1.035984114149125
$ws.Range("J5").Value = 1.040675993702614
$ws.Range("K5").Value = 1.040051129498008
$ws.Range("L5").Value = 1.047829249893798
$ws.Range("M5").Value = 1.059416308752027
$ws.Range("N5").Value = 1.017475999619349

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.03629531326385
$ws.Range("D6").Value = 1.037661262810096
$ws.Range("E6").Value = 1.045468078416783
$ws.Range("F6").Value = 1.057089362912627
$ws.Range("I6").Value = 1.035995682399304
$ws.Range("J6").Value = 1.040712342993529
$ws.Range("K6").Value = 1.040080658472195
$ws.Range("L6").Value = 1.047868455225591
$ws.Range("M6").Value = 1.059461986802008
$ws.Range("N6").Value = 1.017488221503805

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.035932438375445
$ws.Range("D7").Value = 1.037397768272485
$ws.Range("E7").Value = 1.045138444214758
$ws.Range("F7").Value = 1.05671585928185
$ws.Range("I7").Value = 1.035915982243552
$ws.Range("J7").Value = 1.040462333756591
$ws.Range("K7").Value = 1.039877513243993
$ws.Range("L7").Value = 1.047598831893637
$ws.Range("M7").Value = 1.059147864215112
$ws.Range("N7").Value = 1.017404150497859

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.034415120774112
$ws.Range("D8").Value = 1.036295252634752
$ws.Range("E8").Value = 1.043760842685931
$ws.Range("F8").Value = 1.055155046077198
$ws.Range("I8").Value = 1.035578994725245
$ws.Range("J8").Value = 1.039415721368547
$ws.Range("K8").Value = 1.039025945437176
$ws.Range("L8").Value = 1.046470848109808
$ws.Range("M8").Value = 1.057834088563828
$ws.Range("N8").Value = 1.017051973762984

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031739211268269
$ws.Range("D9").Value = 1.034348129616061
$ws.Range("E9").Value = 1.041334103515173
$ws.Range("F9").Value = 1.052406031780915
$ws.Range("I9").Value = 1.034970649330802
$ws.Range("J9").Value = 1.037565331010164
$ws.Range("K9").Value = 1.037516105440131
$ws.Range("L9").Value = 1.044479378384294
$ws.Range("M9").Value = 1.055516003649793
$ws.Range("N9").Value = 1.016428468663685

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.02995370197429
$ws.Range("D10").Value = 1.033047110955577
$ws.Range("E10").Value = 1.039716752819831
$ws.Range("F10").Value = 1.050574194969744
$ws.Range("I10").Value = 1.034555311417743
$ws.Range("J10").Value = 1.036327579067852
$ws.Range("K10").Value = 1.036503293136907
$ws.Range("L10").Value = 1.043149124223431
$ws.Range("M10").Value = 1.053968517244703
$ws.Range("N10").Value = 1.016010823218462

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02918016372438
$ws.Range("D11").Value = 1.032483055181078
$ws.Range("E11").Value = 1.039016524382364
$ws.Range("F11").Value = 1.049781174083657
$ws.Range("I11").Value = 1.034373145446279
$ws.Range("J11").Value = 1.035790622764743
$ws.Range("K11").Value = 1.036063245790699
$ws.Range("L11").Value = 1.042572481084903
$ws.Range("M11").Value = 1.05329792928916
$ws.Range("N11").Value = 1.015829507266097

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02889277504454
$ws.Range("D12").Value = 1.032273433218811
$ws.Range("E12").Value = 1.038756441635557
$ws.Range("F12").Value = 1.049486635983442
$ws.Range("I12").Value = 1.034305131815932
$ws.Range("J12").Value = 1.035591021870554
$ws.Range("K12").Value = 1.035899567734783
$ws.Range("L12").Value = 1.042358193894471
$ws.Range("M12").Value = 1.053048764373597
$ws.Range("N12").Value = 1.015762087186132

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028954423734037
$ws.Range("D13").Value = 1.03231840266944
$ws.Range("E13").Value = 1.038812229683355
$ws.Range("F13").Value = 1.049549814288825
$ws.Range("I13").Value = 1.034319736770866
$ws.Range("J13").Value = 1.035633843807019
$ws.Range("K13").Value = 1.035934687419045
$ws.Range("L13").Value = 1.042404163595342
$ws.Range("M13").Value = 1.053102214673241
$ws.Range("N13").Value = 1.015776552250815

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029156409364911
$ws.Range("D14").Value = 1.032465729940154
$ws.Range("E14").Value = 1.038995025596508
$ws.Range("F14").Value = 1.049756826966906
$ws.Range("I14").Value = 1.03436753054084
$ws.Range("J14").Value = 1.035774126786875
$ws.Range("K14").Value = 1.036049720702383
$ws.Range("L14").Value = 1.042554770007995
$ws.Range("M14").Value = 1.05327733485145
$ws.Range("N14").Value = 1.015823935756387

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029280851088569
$ws.Range("D15").Value = 1.032556489003824
$ws.Range("E15").Value = 1.039107653919945
$ws.Range("F15").Value = 1.049884377604189
$ws.Range("I15").Value = 1.034396931596588
$ws.Range("J15").Value = 1.035860539668329
$ws.Range("K15").Value = 1.036120566688712
$ws.Range("L15").Value = 1.042647550806207
$ws.Range("M15").Value = 1.053385221711159
$ws.Range("N15").Value = 1.01585312084783

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030005030551358
$ws.Range("D16").Value = 1.033084530564723
$ws.Range("E16").Value = 1.039763226592129
$ws.Range("F16").Value = 1.050626828758048
$ws.Range("I16").Value = 1.034567352233212
$ws.Range("J16").Value = 1.036363193907419
$ws.Range("K16").Value = 1.036532466126108
$ws.Range("L16").Value = 1.043187380683798
$ws.Range("M16").Value = 1.05401301098667
$ws.Range("N16").Value = 1.016022846580872

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030459180517072
$ws.Range("D17").Value = 1.033415567726822
$ws.Range("E17").Value = 1.040174474754237
$ws.Range("F17").Value = 1.051092595061461
$ws.Range("I17").Value = 1.034673630788843
$ws.Range("J17").Value = 1.036678226825677
$ws.Range("K17").Value = 1.036790439651962
$ws.Range("L17").Value = 1.043525831209049
$ws.Range("M17").Value = 1.05440666762268
$ws.Range("N17").Value = 1.016129184319308

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030724039950788
$ws.Range("D18").Value = 1.033608588004105
$ws.Range("E18").Value = 1.040414358092281
$ws.Range("F18").Value = 1.051364285791266
$ws.Range("I18").Value = 1.03473539725734
$ws.Range("J18").Value = 1.036861883626236
$ws.Range("K18").Value = 1.036940767237895
$ws.Range("L18").Value = 1.043723182537233
$ws.Range("M18").Value = 1.054636231156671
$ws.Range("N18").Value = 1.016191163728905

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030814343690579
$ws.Range("D19").Value = 1.033674391393331
$ws.Range("E19").Value = 1.040496153709862
$ws.Range("F19").Value = 1.051456928254902
$ws.Range("I19").Value = 1.034756420037198
$ws.Range("J19").Value = 1.03692448949766
$ws.Range("K19").Value = 1.03699200065135
$ws.Range("L19").Value = 1.043790463899199
$ws.Range("M19").Value = 1.054714497990466
$ws.Range("N19").Value = 1.016212289379637

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030410458525387
$ws.Range("D20").Value = 1.033380057619051
$ws.Range("E20").Value = 1.040130350769527
$ws.Range("F20").Value = 1.051042620983287
$ws.Range("I20").Value = 1.034662251278151
$ws.Range("J20").Value = 1.036644436770376
$ws.Range("K20").Value = 1.036762776433437
$ws.Range("L20").Value = 1.043489525001153
$ws.Range("M20").Value = 1.054364437142173
$ws.Range("N20").Value = 1.016117780004041

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.029096931361379
$ws.Range("D21").Value = 1.032422348675445
$ws.Range("E21").Value = 1.038941196397938
$ws.Range("F21").Value = 1.049695866209095
$ws.Range("I21").Value = 1.034353466104985
$ws.Range("J21").Value = 1.035732821133192
$ws.Range("K21").Value = 1.036015852480027
$ws.Range("L21").Value = 1.04251042283802
$ws.Range("M21").Value = 1.053225768505236
$ws.Range("N21").Value = 1.015809984464136

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.028270704279764
$ws.Range("D22").Value = 1.031819582455095
$ws.Range("E22").Value = 1.03819360413121
$ws.Range("F22").Value = 1.048849252956243
$ws.Range("I22").Value = 1.034157300834834
$ws.Range("J22").Value = 1.035158775685642
$ws.Range("K22").Value = 1.035544930374143
$ws.Range("L22").Value = 1.041894265075078
$ws.Range("M22").Value = 1.052509386484384
$ws.Range("N22").Value = 1.015616048769609

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028708737686856
$ws.Range("D23").Value = 1.03213917883171
$ws.Range("E23").Value = 1.038589910046197
$ws.Range("F23").Value = 1.049298045378024
$ws.Range("I23").Value = 1.034261483287893
$ws.Range("J23").Value = 1.035463171439524
$ws.Range("K23").Value = 1.035794698644047
$ws.Range("L23").Value = 1.042220955116886
$ws.Range("M23").Value = 1.052889197743301
$ws.Range("N23").Value = 1.015718896924855

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03043247401482
$ws.Range("D24").Value = 1.033396103318411
$ws.Range("E24").Value = 1.040150288470014
$ws.Range("F24").Value = 1.051065202064759
$ws.Range("I24").Value = 1.034667393881236
$ws.Range("J24").Value = 1.036659705339952
$ws.Range("K24").Value = 1.036775276695845
$ws.Range("L24").Value = 1.043505930402294
$ws.Range("M24").Value = 1.054383519431772
$ws.Range("N24").Value = 1.016122933263998

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032431268119112
$ws.Range("D25").Value = 1.0348520242585
$ws.Range("E25").Value = 1.04196138680138
$ws.Range("F25").Value = 1.053116565483038
$ws.Range("I25").Value = 1.035129643602677
$ws.Range("J25").Value = 1.038044431264214
$ws.Range("K25").Value = 1.037907535743703
$ws.Range("L25").Value = 1.04499467766757
$ws.Range("M25").Value = 1.056115650017009
$ws.Range("N25").Value = 1.016590007183738
